# Actualiza la tabla de empleados: nuevo orden/valores de ingresos y ventas,
# agregando dos empleados nuevos (filas 7-8) y recolocando los dos
# empleados previos al final de la tabla (filas 9-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Writes $text as a genuine text-typed value (not auto-converted to a
    # number) without disturbing the workbook's style table: build the text
    # via a formula (so it comes back as a literal string), copy it, then
    # paste-special *values only* into the destination cell.
    $stage = $ws.Range("ZZ1")
    $escaped = $text -replace '"', '""'
    $stage.Formula = '="' + $escaped + '"'
    $stage.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues
    $stage.Clear()
}

# --- Row 7: IdEmpleado 1 -> 4 (Enrique Soto Andrade), ventas 28 -> 6, ingresos 652.0 -> 128.0
Set-TextCell "A7" "4"
Set-TextCell "B7" "Enrique"
Set-TextCell "C7" "Soto"
Set-TextCell "D7" "Andrade"
$ws.Range("E7").Value = 6.0
Set-TextCell "F7" "128.0"

# --- Row 8: IdEmpleado 3 -> 2 (Aaron Alfonseca Martinez), ventas se mantiene en 5, ingresos 404.0 -> 50.0
Set-TextCell "A8" "2"
Set-TextCell "B8" "Aaron"
Set-TextCell "C8" "Alfonseca"
Set-TextCell "D8" "Martinez"
$ws.Range("E8").Value = 5.0
Set-TextCell "F8" "50.0"

# --- Nueva fila 9: antiguo empleado 1 (Damian Cazarin Montane), ventas 3, ingresos 25.0
$ws.Range("A8:F8").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)  # xlPasteFormats (copia el estilo de la fila 8)

Set-TextCell "A9" "1"
Set-TextCell "B9" "Damian"
Set-TextCell "C9" "Cazarin"
Set-TextCell "D9" "Montane"
$ws.Range("E9").Value = 3.0
Set-TextCell "F9" "25.0"

# --- Nueva fila 10: antiguo empleado 3 (Erick Raymundo Gonzalez Virgen), ventas 2, ingresos 16.0
$ws.Range("A8:F8").Copy()
$ws.Range("A10:F10").PasteSpecial(-4122)  # xlPasteFormats

Set-TextCell "A10" "3"
Set-TextCell "B10" "Erick Raymundo"
Set-TextCell "C10" "Gonzalez"
Set-TextCell "D10" "Virgen"
$ws.Range("E10").Value = 2.0
Set-TextCell "F10" "16.0"

$ws.Range("ZZ1").Clear()
